# Fix Training Data Issue
# Data was taken from 1 day off due to way NBA stats were shown.
# Correct the "Date" column (BF) values from "6-7-2011-12" to "2012-06-07"
# for all data rows (rows 2 through 31).
#
# Note: a plain `Range.Value2 = "2012-06-07"` assignment gets auto-parsed
# as a date by Excel's smart-entry logic (since the text looks like an
# ISO date), which would silently turn the cell into a numeric date
# serial value instead of leaving literal text. To avoid that, we stage
# the literal text in a helper cell that has been explicitly formatted
# as Text, then copy just the *value* (not the format) into each target
# cell with PasteSpecial so the target cells keep their original
# (default) style while still containing the exact text "2012-06-07".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctDate = "2012-06-07"
$oldDate = "6-7-2011-12"

# Helper/staging cell far outside the used data range.
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"
$helper.Value2 = $correctDate

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value2 -eq $oldDate) {
        $helper.Copy()
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}

# Clean up the helper cell so it leaves no residue in the sheet.
$helper.Clear()
$excel.CutCopyMode = $false
